$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 656, shifting existing rows 656:711 down to 657:712
$ws.Rows.Item(656).Insert()

# Populate the newly inserted row 656 with the new record
$ws.Range("A656").Value = 10
$ws.Range("B656").Value = "Vega Modelo de Temuco"
$ws.Range("C656").Value = "La Araucanía"
$ws.Range("D656").Value = 45166
$ws.Range("E656").Value = 9
$ws.Range("F656").Value = 100112008
$ws.Range("G656").Value = "Coliflor"
$ws.Range("H656").Value = "Sin especificar"
$ws.Range("I656").Value = "Primera"
$ws.Range("J656").Value = 1500
$ws.Range("K656").Value = 1300
$ws.Range("L656").Value = 1300
$ws.Range("M656").Value = 1300
$ws.Range("N656").Value = "$/unidad"
$ws.Range("O656").Value = "Provincia del Elquí"
$ws.Range("P656").Value = 1300
$ws.Range("Q656").Value = 1
$ws.Range("R656").Value = "Hortaliza"
